$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old backup code rows A2:A16 so that row 16 (no longer used) is emptied
$ws.Range("A2:A16").ClearContents()

# New backup codes, placed per the target layout:
#  A2:A4 hold three codes, A5:A9 stay blank, A10:A15 hold the remaining six codes
$ws.Range("A2").Value = "95SY88G93C56"
$ws.Range("A3").Value = "HKDJB5BA6J6M"
$ws.Range("A4").Value = "04F5PQ59MWV6"

$ws.Range("A10").Value = "C9DEXVFAR31A"
$ws.Range("A11").Value = "05ANDJ337D9B"
$ws.Range("A12").Value = "SAS5DZQK4GHR"
$ws.Range("A13").Value = "YC7CEVJY9735"
$ws.Range("A14").Value = "X36KP2Z510RZ"
$ws.Range("A15").Value = "MK83F9RSV97N"

# Match the selection shown in the target file
$ws.Range("A5").Select()
